# "modificado por cuarta vez"
# The paragraph currently reads "... modifico por tercera vez el readme".
# We change "tercera" -> "cuarta". Word's own editing history also moves the
# hidden "_GoBack" bookmark (last-edit marker) from right after "tercera" to
# right before the newly-typed "cuarta", which is why the run that used to
# hold " modifico por tercera" ends up split into " modifico por" + " " on
# one side of the bookmark and a fresh "cuarta" run on the other side.

$d = $word.ActiveDocument

# --- Step 1: mark the boundary between "por" and " tercera" -------------
# (temporary bookmark forces Word to keep this as a distinct run boundary)
$rngPorTercera = $d.Content
$rngPorTercera.Find.Execute("por tercera")
$splitAfterPor = $rngPorTercera.Start + 3
$d.Bookmarks.Add("tmpSplitA", $d.Range($splitAfterPor, $splitAfterPor))

# --- Step 2: mark the boundary between " " and "tercera" -----------------
$rngTercera1 = $d.Content
$rngTercera1.Find.Execute("tercera")
$splitBeforeTercera = $rngTercera1.Start
$d.Bookmarks.Add("tmpSplitB", $d.Range($splitBeforeTercera, $splitBeforeTercera))

# --- Step 3: rewrite the now-isolated single space run so it becomes a ---
#     plain, freshly authored run (change then restore, since assigning the
#     identical text is otherwise a no-op for run authorship purposes).
$spaceRng = $d.Range($splitAfterPor, $splitBeforeTercera)
$spaceRng.Text = "Z"
$spaceRng2 = $d.Range($splitAfterPor, $splitAfterPor + 1)
$spaceRng2.Text = " "

# --- Step 4: replace "tercera" with "cuarta" (still isolated between -----
#     tmpSplitB and the original _GoBack bookmark, so it becomes its own
#     cleanly authored run without inheriting neighbouring whitespace).
$rngTercera2 = $d.Content
$rngTercera2.Find.Execute("tercera")
$rngTercera2.Text = "cuarta"

# --- Step 5: move "_GoBack" so it again sits right before the edited word,
#     i.e. exactly where tmpSplitB currently is. Re-adding the reserved
#     "_GoBack" name relocates Word's single hidden bookmark automatically.
$rngCuarta = $d.Content
$rngCuarta.Find.Execute("cuarta")
$d.Bookmarks.Add("_GoBack", $d.Range($rngCuarta.Start, $rngCuarta.Start))

# --- Step 6: drop the scaffolding bookmarks; the run splits they created -
#     remain in place even after the bookmarks themselves are gone.
$d.Bookmarks.Item("tmpSplitA").Delete()
$d.Bookmarks.Item("tmpSplitB").Delete()
